# Applies the "Lake Florida sample" update: appends 25 new survey response
# rows (rows 348-372) to the "Form Responses 1" sheet, reproducing the data
# captured by the linked Google Form between 2018-08-02 and 2018-08-16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A template data row (row 2) supplies the cell formatting that each new
# row should inherit: column A keeps the timestamp datetime format, column
# B keeps the date format, and every other column keeps the plain/general
# format already used throughout the sheet.
$FormatTemplateRow = 2
$ColumnIndex = @{ A=1; B=2; C=3; D=4; E=5; F=6; G=7; H=8; I=9; J=10; K=11; L=12; M=13 }

function Copy-CellFormat($targetRow, $targetCol, $templateCol) {
    $ws.Cells.Item($FormatTemplateRow, $templateCol).Copy() | Out-Null
    $ws.Cells.Item($targetRow, $targetCol).PasteSpecial(-4122) | Out-Null
}

$NewRows = @(
    @{ Row=348; A=43314.50230461806; B=43313.0; C="Lake Burgan"; D="Austen"; E=11.0; F=14.1; G=82.0; H="Right"; I=1.0; J=1.5; K=1.0; L="Rock" },
    @{ Row=349; A=43314.50268737269; B=43313.0; C="Lake Burgan"; D="Austen"; E=11.0; F=14.52; G=60.0; H="Left"; I=1.0; J=1.2; K=0.6; L="Gravel" },
    @{ Row=350; A=43314.50307409722; B=43313.0; C="Lake Burgan"; D="Austen"; E=11.0; F=15.5; G=45.0; H="Left"; I=1.0; J=3.0; K=1.8; L="Rock" },
    @{ Row=351; A=43314.50361134259; B=43313.0; C="Lake Burgan"; D="Austen"; E=11.0; F=19.41; G=54.0; H="Right"; I=1.0; J=1.5; K=0.9; L="Shellfish"; M="on dead shell" },
    @{ Row=352; A=43314.506624629634; B=43313.0; C="Lake Burgan"; D="Austen"; E=12.0; M="No mussels" },
    @{ Row=353; A=43314.50792142361; B=43313.0; C="Lake Burgan"; D="Austen"; E=13.0; F=3.45; G=17.0; H="Right"; I=1.0; J=1.6; K=1.0; L="Wood" },
    @{ Row=354; A=43314.50824556713; B=43313.0; C="Lake Burgan"; D="Austen"; E=13.0; F=5.31; G=12.0; H="Right"; I=1.0; J=1.6; K=1.0; L="Rock" },
    @{ Row=355; A=43314.508628125; B=43313.0; C="Lake Burgan"; D="Austen"; E=13.0; F=13.18; G=28.0; H="Left"; I=1.0; J=1.6; K=1.0; L="Rock" },
    @{ Row=356; A=43314.50901388889; B=43313.0; C="Lake Burgan"; D="Austen"; E=14.0; F=10.6; G=90.0; H="Left"; I=1.0; J=1.8; K=1.0; L="Rock" },
    @{ Row=357; A=43314.51064493056; B=43313.0; C="Lake Burgan"; D="Austen"; E=15.0; F=3.45; G=41.0; H="Right"; I=2.0; J=5.0; K=2.0; L="Rock" },
    @{ Row=358; A=43314.51115289352; B=43313.0; C="Lake Burgan"; D="Austen"; E=15.0; F=5.12; G=55.0; H="Left"; I=2.0; J=1.8; K=2.0; L="Rock" },
    @{ Row=359; A=43314.51164375; B=43313.0; C="Lake Burgan"; D="Austen"; E=15.0; F=8.02; G=98.0; H="Right"; I=1.0; J=1.7; K=1.0; L="Rock" },
    @{ Row=360; A=43314.5122212037; B=43313.0; C="Lake Burgan"; D="Austen"; E=15.0; F=9.1; G=13.0; H="Left"; I=3.0; J=16.0; K=8.0; L="Wood" },
    @{ Row=361; A=43314.512673136574; B=43313.0; C="Lake Burgan"; D="Austen"; E=15.0; F=9.4; G=14.0; H="Right"; I=1.0; J=1.5; K=0.9; L="Wood" },
    @{ Row=362; A=43314.51318565972; B=43313.0; C="Lake Burgan"; D="Austen"; E=15.0; F=15.45; G=15.0; H="Left"; I=1.0; J=2.0; K=1.2; L="Rock" },
    @{ Row=363; A=43314.51356288194; B=43313.0; C="Lake Burgan"; D="Austen"; E=15.0; F=16.85; G=4.0; H="Left"; I=2.0; J=13.0; K=3.0; L="Rock" },
    @{ Row=364; A=43314.5139537037; B=43313.0; C="Lake Burgan"; D="Austen"; E=15.0; F=17.08; G=93.0; H="Left"; I=1.0; J=2.5; K=1.7; L="Gravel" },
    @{ Row=365; A=43314.514324143514; B=43313.0; C="Lake Burgan"; D="Austen"; E=15.0; F=18.1; G=85.0; H="Left"; I=1.0; J=1.8; K=1.2; L="Rock" },
    @{ Row=366; A=43314.514644826384; B=43313.0; C="Lake Burgan"; D="Austen"; E=15.0; F=19.18; G=25.0; H="Left"; I=1.0; J=2.0; K=1.2; L="Rock" },
    @{ Row=367; A=43328.49418869213; B=43326.0; C="Lake Florida"; D="Austen"; E=13.0; M="0 mussels" },
    @{ Row=368; A=43328.49526739583; B=43326.0; C="Lake Florida"; D="Austen"; E=14.0; M="0 mussels" },
    @{ Row=369; A=43328.4956521875; B=43326.0; C="Lake Florida"; D="Austen"; E=15.0; M="0 mussels" },
    @{ Row=370; A=43328.54446502315; B=43326.0; C="Lake Florida"; D="Aislyn"; E=13.0; I=0.0 },
    @{ Row=371; A=43328.54477920139; B=43326.0; C="Lake Florida"; D="Aislyn"; E=14.0; F=3.27; G=22.0; H="Left"; I=1.0; J=0.4; K=0.2; L="Shellfish" },
    @{ Row=372; A=43328.54500013889; B=43326.0; C="Lake Florida"; D="Aislyn"; E=15.0; I=0.0 },
)

foreach ($rowData in $NewRows) {
    $r = $rowData.Row
    foreach ($col in @("A","B","C","D","E","F","G","H","I","J","K","L","M")) {
        if (-not $rowData.ContainsKey($col)) { continue }

        $colIdx = $ColumnIndex[$col]
        if ($col -eq "A") {
            $templateCol = 1
        } elseif ($col -eq "B") {
            $templateCol = 2
        } else {
            $templateCol = 3
        }

        Copy-CellFormat $r $colIdx $templateCol
        $ws.Cells.Item($r, $colIdx).Value = $rowData[$col]
    }
}
